$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (26.1666... chars -> stored OOXML width of 27)
$ws.Columns.Item(1).ColumnWidth = 26.16666666666667

# Fill email values into A2:A5
$ws.Range("A2").Value = "emily.johnson@example.com"
$ws.Range("A3").Value = "emily.johnson@example.com"
$ws.Range("A4").Value = "emily.johnson@example.com"
$ws.Range("A5").Value = "emily.johnson@example.com"
